# "work on the robot 3/29"
# Update motor-mapping calculations on Sheet1 and clean up scratch cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the scratch column J calculations near the top of the sheet (row 13 & 15)
$ws.Range("J13").ClearContents() | Out-Null
$ws.Range("I15").ClearContents() | Out-Null
$ws.Range("J15").ClearContents() | Out-Null

# Tilt motor block (rows 17-19): update gear ratio input and add new ratio helper cells
$ws.Range("D17").Value = 18
$ws.Range("I17").Value = 66
$ws.Range("I18").Value = 27.5
$ws.Range("I19").Formula = "=I18/I17"

# Up/Down motor block (rows 20-21): update gear ratio input and add new helper cells
$ws.Range("D20").Value = 49.5
$ws.Range("I20").Value = 0.056
$ws.Range("I21").Formula = "=I19*I20"

# Hooks motor block (row 24): remove the extra 5x multiplier formula, keep the cell's style
$ws.Range("H24").ClearContents() | Out-Null

# Column H is no longer used for a long formula result; shrink it back down
$ws.Range("H1").EntireColumn.ColumnWidth = 8.33

# Restore the view to where the user was last working
$ws.Range("G21").Select() | Out-Null
